# Link to Python no-PIE bug.
#
# Slide 15 ("Position-independent Executable (PIE)"), the bullet that
# currently just reads "/usr/bin/python3" becomes
# "/usr/bin/python3 (Launchpad #1452115)" with "Launchpad #1452115"
# hyperlinked to the Launchpad bug tracking Ubuntu's python3 package not
# being built as a PIE.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(2)

$fullTr = $shp.TextFrame.TextRange
$tr = $fullTr.Paragraphs(8, 1)

# Append " (" to the existing run instead of creating a new one, so the
# run boundary matches the target markup exactly.
$run1 = $tr.Runs(1)
$newPrefix = $run1.Text.TrimEnd("`r") + " ("
$run1.Text = $newPrefix

$linkLabel = "Launchpad #1452115"

# Remember where the paragraph (and therefore the new run) starts so we
# can re-select exactly the inserted label afterwards. (Note: reading
# $run1.Text back here would include a trailing "`r" paragraph mark, so
# use the literal string we just assigned instead.)
$paraStart = $tr.Start
$prefixLength = $newPrefix.Length

$null = $tr.InsertAfter($linkLabel)
$null = $tr.InsertAfter(")")

$linkRange = $fullTr.Characters($paraStart + $prefixLength, $linkLabel.Length)
$linkRange.ActionSettings(1).Hyperlink.Address = "https://bugs.launchpad.net/ubuntu/+source/python3.4/+bug/1452115"

Write-Output $tr.Text
